$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 63, pushing existing rows 63..102 down to 64..103
$ws.Rows.Item(63).Insert()

# Populate the newly inserted row 63 with the new weekly data point
$ws.Cells.Item(63, 1).Value = 2
$ws.Cells.Item(63, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(63, 3).Value = "Coquimbo"
$ws.Cells.Item(63, 4).Value = 45141
$ws.Cells.Item(63, 5).Value = 4
$ws.Cells.Item(63, 6).Value = 100112026
$ws.Cells.Item(63, 7).Value = "Haba"
$ws.Cells.Item(63, 8).Value = "Sin especificar"
$ws.Cells.Item(63, 9).Value = "Primera"
$ws.Cells.Item(63, 10).Value = 1400
$ws.Cells.Item(63, 11).Value = 7000
$ws.Cells.Item(63, 12).Value = 8000
$ws.Cells.Item(63, 13).Value = 7500
$ws.Cells.Item(63, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(63, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(63, 16).Value = 300
$ws.Cells.Item(63, 17).Value = 25
$ws.Cells.Item(63, 18).Value = "Hortaliza"
